$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was reported for "Vega Modelo de Temuco" (Espárragos).
# It belongs right after the header block of existing records, at row 52 -
# insert a new row there (this pushes the former rows 52-124 down to 53-125,
# which is exactly what the diff shows: each old row's data now lives one
# row further down, and the table grows from A1:R124 to A1:R125).
$ws.Rows.Item(52).Insert()

# Populate the newly inserted row 52 with the new weekly observation.
$ws.Range("A52").Value = 10
$ws.Range("B52").Value = "Vega Modelo de Temuco"
$ws.Range("C52").Value = "La Araucanía"
$ws.Range("D52").Value = "2023-12-11"
$ws.Range("E52").Value = 9
$ws.Range("F52").Value = 300000000
$ws.Range("G52").Value = "Espárragos"
$ws.Range("H52").Value = "Sin especificar"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 300
$ws.Range("K52").Value = 1800
$ws.Range("L52").Value = 1800
$ws.Range("M52").Value = 1800
$ws.Range("N52").Value = "$/kilo"
$ws.Range("O52").Value = "Región del Maule"
$ws.Range("P52").Value = 1800
$ws.Range("Q52").Value = 1
$ws.Range("R52").Value = "Hortaliza"
